$d = $word.ActiveDocument

$replacements = @(
    @{ Table = 1; Row = 1; Col = 1; Text = 'dame' }
    @{ Table = 1; Row = 1; Col = 2; Text = 'only' }
    @{ Table = 1; Row = 1; Col = 3; Text = 'sweet' }
    @{ Table = 1; Row = 1; Col = 4; Text = 'heart' }
    @{ Table = 1; Row = 1; Col = 5; Text = 'gusta' }
    @{ Table = 1; Row = 2; Col = 1; Text = 'amarte' }
    @{ Table = 1; Row = 2; Col = 2; Text = 'quiere' }
    @{ Table = 1; Row = 2; Col = 3; Text = 'para' }
    @{ Table = 1; Row = 2; Col = 4; Text = 'besos' }
    @{ Table = 1; Row = 2; Col = 5; Text = 'cuando' }
    @{ Table = 1; Row = 3; Col = 1; Text = 'mejor' }
    @{ Table = 1; Row = 3; Col = 2; Text = 'it''s' }
    @{ Table = 1; Row = 3; Col = 3; Text = 'quien' }
    @{ Table = 1; Row = 3; Col = 4; Text = 'hasta' }
    @{ Table = 1; Row = 3; Col = 5; Text = 'vamos' }
    @{ Table = 1; Row = 4; Col = 1; Text = 'boom' }
    @{ Table = 1; Row = 4; Col = 2; Text = 'solo' }
    @{ Table = 1; Row = 4; Col = 4; Text = 'mujer' }
    @{ Table = 1; Row = 4; Col = 5; Text = 'este' }
    @{ Table = 1; Row = 5; Col = 1; Text = 'juro' }
    @{ Table = 1; Row = 5; Col = 2; Text = 'algo' }
    @{ Table = 1; Row = 5; Col = 3; Text = 'llevar' }
    @{ Table = 1; Row = 5; Col = 4; Text = 'dímelo' }
    @{ Table = 1; Row = 5; Col = 5; Text = 'bien' }
    @{ Table = 1; Row = 6; Col = 1; Text = 'dame' }
    @{ Table = 1; Row = 6; Col = 2; Text = 'only' }
    @{ Table = 1; Row = 6; Col = 3; Text = 'sweet' }
    @{ Table = 1; Row = 6; Col = 4; Text = 'heart' }
    @{ Table = 1; Row = 6; Col = 5; Text = 'gusta' }
    @{ Table = 2; Row = 1; Col = 1; Text = 'dura' }
    @{ Table = 2; Row = 1; Col = 2; Text = 'quiere' }
    @{ Table = 2; Row = 1; Col = 3; Text = 'nosotros' }
    @{ Table = 2; Row = 1; Col = 4; Text = 'they' }
    @{ Table = 2; Row = 1; Col = 5; Text = 'digo' }
    @{ Table = 2; Row = 2; Col = 1; Text = 'down' }
    @{ Table = 2; Row = 2; Col = 2; Text = 'quiera' }
    @{ Table = 2; Row = 2; Col = 3; Text = 'llegue' }
    @{ Table = 2; Row = 2; Col = 4; Text = 'algo' }
    @{ Table = 2; Row = 2; Col = 5; Text = 'perdí' }
    @{ Table = 2; Row = 3; Col = 1; Text = 'mami' }
    @{ Table = 2; Row = 3; Col = 2; Text = 'don''t' }
    @{ Table = 2; Row = 3; Col = 3; Text = 'bien' }
    @{ Table = 2; Row = 3; Col = 4; Text = 'luian' }
    @{ Table = 2; Row = 3; Col = 5; Text = 'falta' }
    @{ Table = 2; Row = 4; Col = 1; Text = 'you''re' }
    @{ Table = 2; Row = 4; Col = 2; Text = 'that''s' }
    @{ Table = 2; Row = 4; Col = 4; Text = 'girl' }
    @{ Table = 2; Row = 4; Col = 5; Text = 'fácil' }
    @{ Table = 2; Row = 5; Col = 1; Text = 'it''s' }
    @{ Table = 2; Row = 5; Col = 2; Text = 'solo' }
    @{ Table = 2; Row = 5; Col = 3; Text = 'make' }
    @{ Table = 2; Row = 5; Col = 4; Text = 'pare' }
    @{ Table = 2; Row = 5; Col = 5; Text = 'ganas' }
    @{ Table = 2; Row = 6; Col = 1; Text = 'dura' }
    @{ Table = 2; Row = 6; Col = 2; Text = 'quiere' }
    @{ Table = 2; Row = 6; Col = 3; Text = 'nosotros' }
    @{ Table = 2; Row = 6; Col = 4; Text = 'they' }
    @{ Table = 2; Row = 6; Col = 5; Text = 'digo' }
)

foreach ($item in $replacements) {
    $cell = $d.Tables.Item($item.Table).Cell($item.Row, $item.Col)
    $full = $cell.Range
    # The cell content is: <break>(1 char) <text> <break>(1 char) <para+cell mark>(last unit)
    # Target just the text portion (between the two breaks) so the break runs are preserved.
    $innerStart = $full.Start + 1
    $innerEnd = $full.End - 2
    $target = $d.Range($innerStart, $innerEnd)
    $target.Text = $item.Text
}

Write-Host "Done applying" $replacements.Count "replacements"
